# Apply weekly Fruta/Hortaliza update to "Feria Lagunitas de Puerto Montt - Nectarín" sheet.
# Inserts 4 new price rows (a new week's data, 2023-01-13, 14-kilo boxes) above the
# existing block of 6 rows (which covers 2021-01-13 / 15-kilo boxes), shifting the
# existing rows down by 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows before row 546; this pushes the old rows 546-551 down to 550-555.
$ws.Range("A546:A549").EntireRow.Insert()

# Common values shared by all rows in this block.
$mercadoId = 4
$feria = "Feria Lagunitas de Puerto Montt"
$region = "Los Lagos"
$periodo = 10
$tipo = "Fruta"
$grupoId = 100103
$grupo = "Frutos de hueso (carozo)"
$especieId = 100103006
$especie = "Nectarín"
$regionOrigen = "Región de O'Higgins"

function Set-Row($r, $variedad, $calidad, $cantidad, $precioMin, $precioMax, $precioFrec, $unidad, $precioKilo, $kilos) {
    $ws.Cells.Item($r, 1).Value = $mercadoId
    $ws.Cells.Item($r, 2).Value = $feria
    $ws.Cells.Item($r, 3).Value = $region
    $ws.Cells.Item($r, 5).Value = $periodo
    $ws.Cells.Item($r, 6).Value = $tipo
    $ws.Cells.Item($r, 7).Value = $grupoId
    $ws.Cells.Item($r, 8).Value = $grupo
    $ws.Cells.Item($r, 9).Value = $especieId
    $ws.Cells.Item($r, 10).Value = $especie
    $ws.Cells.Item($r, 11).Value = $variedad
    $ws.Cells.Item($r, 12).Value = $calidad
    $ws.Cells.Item($r, 13).Value = $cantidad
    $ws.Cells.Item($r, 14).Value = $precioMin
    $ws.Cells.Item($r, 15).Value = $precioMax
    $ws.Cells.Item($r, 16).Value = $precioFrec
    $ws.Cells.Item($r, 17).Value = $unidad
    $ws.Cells.Item($r, 18).Value = $regionOrigen
    $ws.Cells.Item($r, 19).Value = $precioKilo
    $ws.Cells.Item($r, 20).Value = $kilos
}

# New date for the inserted week's rows (2023-01-13).
$ws.Range("D546:D549").Value = 44939

Set-Row 546 "Super Queen" "Especial" 300 23000 23000 23000 "`$/caja 14 kilos empedrada" 1643 14
Set-Row 547 "Super Queen" "Primera"  600 18000 19000 18500 "`$/caja 14 kilos empedrada" 1321 14
Set-Row 548 "Venus"       "Especial" 300 23000 23000 23000 "`$/caja 14 kilos empedrada" 1643 14
Set-Row 549 "Venus"       "Primera"  600 18000 19000 18500 "`$/caja 14 kilos empedrada" 1321 14

Write-Host "done"
